# Edit script: add latest alkalinity standard run data (04/07/2021 and 04/27/2021
# batches) to the main data sheet, plus the corresponding BayStd2 summary row on
# the "QAQC baystds assessment" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data 11dec2019 to 17nov2020": append new sample rows 190-206
# ---------------------------------------------------------------------------
$data = $wb.Worksheets.Item("data 11dec2019 to 17nov2020")

# Row 190: SAC-2
$data.Range("A190").Value = "SAC-2"
$data.Range("B190").Value = 2290.3416417366602

# Row 191: BAYSTD1-04072021
$data.Range("A191").Value = "BAYSTD1-04072021"
$data.Range("B191").Value = 2185.6742415221702

# Row 192: SAC-1
$data.Range("A192").Value = "SAC-1"
$data.Range("B192").Value = 2450.9836064696501

# Row 193-194: P-0084-1 / P-0084-2 (with stats on row 194)
$data.Range("A193").Value = "P-0084-1"
$data.Range("B193").Value = 2203.76908601887
$data.Range("A194").Value = "P-0084-2"
$data.Range("B194").Value = 2339.16249328406
$data.Range("E194").Formula = "=AVERAGE(B193:B194)"
$data.Range("F194").Formula = "=STDEV.S(B193:B194)"
$data.Range("G194").Formula = "=2*F194"
$data.Range("H194").Formula = "=F194/E194"
$data.Range("I194").Formula = "=H194"
$data.Range("J194").Formula = "=MIN(B193:B194)"
$data.Range("K194").Formula = "=MAX(B193:B194)"
$data.Range("L194").Formula = "=K194-J194"
$data.Range("E194:L194").Style = "QAQC_style"
$data.Range("I194").Style = "QAQC_pct_style"

# Row 195: SAC-2
$data.Range("A195").Value = "SAC-2"
$data.Range("B195").Value = 2380.6235945409098

# Row 196: DIC-CRM186-04272021-1 (excluded, note in D)
$data.Range("A196").Value = "DIC-CRM186-04272021-1"
$data.Range("B196").Value = 2250.7575623355901
$data.Range("D196").Value = "not sure what happened here, might be worth excluding"

# Row 197-198: DIC-CRM186-04272021-2/3 (stats use only rows 197:198)
$data.Range("A197").Value = "DIC-CRM186-04272021-2"
$data.Range("B197").Value = 2203.22618654695
$data.Range("A198").Value = "DIC-CRM186-04272021-3"
$data.Range("B198").Value = 2200.2796088885002
$data.Range("D198").Value = "using only two values in the calulations to the right"
$data.Range("E198").Formula = "=AVERAGE(B197:B198)"
$data.Range("F198").Formula = "=STDEV.S(B197:B198)"
$data.Range("G198").Formula = "=2*F198"
$data.Range("H198").Formula = "=F198/E198"
$data.Range("I198").Formula = "=H198"
$data.Range("J198").Formula = "=MIN(B197:B198)"
$data.Range("K198").Formula = "=MAX(B197:B198)"
$data.Range("L198").Formula = "=K198-J198"

# Row 199: SAC-1
$data.Range("A199").Value = "SAC-1"
$data.Range("B199").Value = 2398.0314183502901

# Rows 200-202: BAYSTD1-04272021-1/2/3 (stats on row 202)
$data.Range("A200").Value = "BAYSTD1-04272021-1"
$data.Range("B200").Value = 2188.0548117417202
$data.Range("A201").Value = "BAYSTD1-04272021-2"
$data.Range("B201").Value = 2182.1633268279402
$data.Range("A202").Value = "BAYSTD1-04272021-3"
$data.Range("B202").Value = 2179.9038352182201
$data.Range("E202").Formula = "=AVERAGE(B200:B202)"
$data.Range("F202").Formula = "=STDEV.S(B200:B202)"
$data.Range("G202").Formula = "=2*F202"
$data.Range("H202").Formula = "=F202/E202"
$data.Range("I202").Formula = "=H202"
$data.Range("J202").Formula = "=MIN(B200:B202)"
$data.Range("K202").Formula = "=MAX(B200:B202)"
$data.Range("L202").Formula = "=K202-J202"

# Row 203: SAC-3
$data.Range("A203").Value = "SAC-3"
$data.Range("B203").Value = 2110.5869142463098

# Rows 204-206: BAYSTD2-04272021-1/2/3 (stats on row 206)
$data.Range("A204").Value = "BAYSTD2-04272021-1"
$data.Range("B204").Value = 2253.49247123941
$data.Range("A205").Value = "BAYSTD2-04272021-2"
$data.Range("B205").Value = 2256.3586638595498
$data.Range("A206").Value = "BAYSTD2-04272021-3"
$data.Range("B206").Value = 2252.4539135216601
$data.Range("E206").Formula = "=AVERAGE(B204:B206)"
$data.Range("F206").Formula = "=STDEV.S(B204:B206)"
$data.Range("G206").Formula = "=2*F206"
$data.Range("H206").Formula = "=F206/E206"
$data.Range("I206").Formula = "=H206"
$data.Range("J206").Formula = "=MIN(B204:B206)"
$data.Range("K206").Formula = "=MAX(B204:B206)"
$data.Range("L206").Formula = "=K206-J206"

# ---------------------------------------------------------------------------
# Sheet "QAQC baystds assessment": add BayStd2 summary reference + new row 29
# ---------------------------------------------------------------------------
$baystds = $wb.Worksheets.Item("QAQC baystds assessment")

$baystds.Range("K6").Value = "BayStd2"
$baystds.Range("L6").Value = 2254.1016828735396
$baystds.Range("M6").Value = 30.34

$baystds.Range("A29").Value = "BAYSTD1-04072021"
$baystds.Range("B29").Value = 2185.6742415221702
$baystds.Range("C29").Formula = "=`$L`$5"
$baystds.Range("D29").Formula = "=B29-C29"
$baystds.Range("E29").Formula = "=D29^2"
